$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.140.20"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "'3.275.27"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'588.12"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").Value = "'186.57"
$ws.Range("E6").Value = "  +3.92%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.601"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +4.28%  "
$ws.Range("D10").Value = "'6.74"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").Value = "'0.417"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "'3.844.53"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").Value = "'28.67"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").Value = "'68.148.91"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("E16").Value = "  +2.46%  "
$ws.Range("D17").Value = "'3.277.07"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "'13.63"
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").Value = "'7.73"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("D26").Value = "'0.190"
$ws.Range("E26").Value = "  +6.04%  "
$ws.Range("D27").Value = "'9.79"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'5.81"
$ws.Range("E29").Value = "  +3.12%  "
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("D31").Value = "'22.94"
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D32").Value = "'7.16"
$ws.Range("E32").Value = "  +5.17%  "
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("E35").Value = "  +2.47%  "
$ws.Range("D36").Value = "'163.14"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").Value = "'1.87"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("E39").Value = "  +2.45%  "
$ws.Range("D40").Value = "'26.67"
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("D41").Value = "'4.62"
$ws.Range("E41").Value = "  +4.64%  "
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("D43").Value = "'41.36"
$ws.Range("E43").Value = "  +2.16%  "
$ws.Range("D45").Value = "'25.47"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("D46").Value = "'345.15"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "'2.644.04"
$ws.Range("E47").Value = "  -4.58%  "
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("D49").Value = "'32.08"
$ws.Range("E49").Value = "  +3.73%  "
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("E51").Value = "  +0.01%  "
